$wb = $excel.ActiveWorkbook

# --- Sheet "Add Devices" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Add Devices")
$ws1.Range("K8").Value = "IOB800(x1)"
$ws1.Range("K9").Value = "AttachedFunctionality"
$ws1.Range("O10").Value = "Other Slot Cards  (3 of 18"

# --- Sheet "Defect 1559" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Defect 1559")
$ws2.Range("E8").Value = "IOB800(x1)"
$ws2.Range("E9").Value = "AttachedFunctionality"
$ws2.Range("G10").Value = 1
$ws2.Range("I10").Value = 0
$ws2.Range("O10").Value = "Other Slot Cards  (7 of 17"

# --- Sheet "Defect 1545" (sheet3) ---
$ws3 = $wb.Worksheets.Item("Defect 1545")
$ws3.Range("D8").Value = "IOB800(x2)"
$ws3.Range("D9").Value = "AttachedFunctionality"
$ws3.Range("D10").Value = 1
$ws3.Range("N10").Value = "Other Slot Cards  (3 of 18"

# --- Sheet "Defect 1545(2)" (sheet4) ---
$ws4 = $wb.Worksheets.Item("Defect 1545(2)")
$ws4.Range("G10").Value = 1
$ws4.Range("I10").Value = 0
$ws4.Range("N10").Value = "Other Slot Cards  (2 of 2"

# --- Selections ---
$ws1.Range("E8").Select()
$ws2.Range("E8:E9").Select()
$ws3.Range("M9").Select()

# Activate "Defect 1545(2)" last so it becomes the active/selected tab,
# and select L9 on it.
$ws4.Activate()
$ws4.Range("L9").Select()
